$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 14: "Want something typically Japanese?" question block ---
$ws.Range("A14").Value = 3
$ws.Range("B14").Value = "Want something typically Japanese?"
$ws.Range("C14").Value = 1
$ws.Range("D14").Value = "Samurai are cool (even cooler with hip-hop music)"
$ws.Range("G14").Value = "Samurai Champloo"

$ws.Range("C15").Value = 2
$ws.Range("D15").Value = "I'VE HEARD THERE ARE SOME GIANT ROBOTS"
$ws.Range("G15").Value = "TENGEN TOPPA GURREN LAGANN"

$ws.Range("C16").Value = 3
$ws.Range("D16").Value = "I like the Ghibli movies, especially the spirit stuff"
$ws.Range("G16").Value = "Mushishi"

$ws.Range("C17").Value = 4
$ws.Range("D17").Value = "Something more Western, please"
$ws.Range("F17").Value = 4

# --- Row 18: "Fantasy action adventure?" question block ---
$ws.Range("A18").Value = 4
$ws.Range("B18").Value = "Fantasy action adventure?"
$ws.Range("C18").Value = 1
$ws.Range("D18").Value = "Something gritty, edgy and with epic music"
$ws.Range("G18").Value = "Attack on Titan"

$ws.Range("C19").Value = 2
$ws.Range("D19").Value = "Arabian Nights road trip adventure, Avatar-style"
$ws.Range("G19").Value = "Magi: The Labyrinth of Magic"

$ws.Range("C20").Value = 3
$ws.Range("D20").Value = "Pirates of the Carribean meets Diablo"
$ws.Range("G20").Value = "Shingeki no Bahamut: Genesis"

$ws.Range("C21").Value = 4
$ws.Range("D21").Value = "I want more!"
$ws.Range("G21").Value = "Extended version in development!"

# --- Update existing rows 10-13 (add "Next Question" labels in G, and F13 pointer) ---
$ws.Range("G10").Value = "Psycho-Pass"
$ws.Range("G11").Value = "Darker than Black"
$ws.Range("G12").Value = "Steins;Gate"
$ws.Range("F13").Value = 3

# --- Extend the shared CONCATENATE formula down through E35 ---
# (kept as a separate range from the pre-existing E3:E13 shared group so the
#  engine doesn't corrupt the existing group's master-cell formula text)
$ws.Range("E14:E35").Formula = '=CONCATENATE("[",C14,"] ",D14)'

# --- Row heights: taller (30) rows for the longer two-line choice text ---
$ws.Rows.Item(14).RowHeight = 30
$ws.Rows.Item(15).RowHeight = 30
$ws.Rows.Item(20).RowHeight = 30
$ws.Rows.Item(21).RowHeight = 30

# --- Column widths: D grows to fit "Something gritty..." text, E shrinks ---
$ws.Columns.Item(4).ColumnWidth = 22.5
$ws.Columns.Item(5).ColumnWidth = 12.333333333333334

# --- Scroll back to top-left and move the active selection to A18 ---
$ws.Application.ActiveWindow.ScrollColumn = 1
$ws.Application.ActiveWindow.ScrollRow = 1
$ws.Range("A18").Select()
